# Commit: "Fruta / hortaliza, semanal" — weekly refresh of the Espinaca
# (Mercado Mayorista Lo Valledor de Santiago) dataset. A brand-new weekly
# observation is inserted at row 628, pushing the previously-existing rows
# 628-667 down to 629-668 (their contents are unchanged, only their row
# numbers shift). The sheet's used range grows from A1:R667 to A1:R668.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 628, shifting rows
# 628:667 down to 629:668 (all of their data travels with them).
$ws.Rows.Item(628).Insert()

# Populate the newly inserted row 628 with the new weekly record.
$ws.Cells.Item(628, 1).Value = 6
$ws.Cells.Item(628, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(628, 3).Value = "Metropolitana"
$ws.Cells.Item(628, 4).Value = 44931
$ws.Cells.Item(628, 5).Value = 13
$ws.Cells.Item(628, 6).Value = 100112012
$ws.Cells.Item(628, 7).Value = "Espinaca"
$ws.Cells.Item(628, 8).Value = "Sin especificar"
$ws.Cells.Item(628, 9).Value = "Primera"
$ws.Cells.Item(628, 10).Value = 450
$ws.Cells.Item(628, 11).Value = 7000
$ws.Cells.Item(628, 12).Value = 8000
$ws.Cells.Item(628, 13).Value = 7422
$ws.Cells.Item(628, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(628, 15).Value = "Región Metropolitana"
$ws.Cells.Item(628, 16).Value = 742
$ws.Cells.Item(628, 17).Value = 10
$ws.Cells.Item(628, 18).Value = "Hortaliza"
